$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$overview.Range("G3").Value = "2016-08-12 12:53:22"
$zhcn.Range("H3").Value     = "2016-08-12 12:53:15"
$zhcn.Range("K3").Value     = "2016-08-12 12:53:44"
$dede.Range("H3").Value     = "2016-08-12 12:53:22"
$dede.Range("K3").Value     = "2016-08-12 12:53:54"
